# New weekly price observation arrived for
# "Hortaliza, Terminal Hortofrutícola Agro Chillán - Arveja Verde".
# The data block lives in rows 9-31 (most recent first). A new row is
# inserted at the top (row 9) for the latest reading, every existing
# row shifts down by one, and the oldest reading (which was in row 31)
# drops out of the window.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 9-31 down by one: insert a fresh row 9 (carries D9's date
# style along), then drop the row that is now 32 (the old row 31 that
# fell out of the trailing window) so rows 32+ are undisturbed.
$ws.Rows.Item(9).Insert()
$ws.Rows.Item(32).Delete()

# Populate the newly inserted row 9 with the latest observation.
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C9").Value = "Ñuble"
$ws.Range("D9").Value = 44565
$ws.Range("E9").Value = 16
$ws.Range("F9").Value = 100112022
$ws.Range("G9").Value = "Arveja Verde"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 18000
$ws.Range("L9").Value = 19000
$ws.Range("M9").Value = 18500
$ws.Range("N9").Value = "$/saco 25 kilos"
$ws.Range("O9").Value = "Provincia de Diguillín"
$ws.Range("P9").Value = 740
$ws.Range("Q9").Value = 25
$ws.Range("R9").Value = "Hortaliza"
